$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.804.86'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '3.366.30'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.87%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.68%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.64'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.73%  '
$ws.Range("E10").Value = '  -2.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.381'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.68%  '
$ws.Range("D12").Value = '3.941.25'
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("E13").Value = '  +1.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.33%  '
$ws.Range("D15").Value = '3.356.32'
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("E16").Value = '  -1.99%  '
$ws.Range("D17").Value = '60.925.74'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.87'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '380.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.547'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.24%  '
$ws.Range("E24").Value = '  +0.27%  '
$ws.Range("E25").Value = '  -5.33%  '
$ws.Range("E26").Value = '  +6.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.53%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.82'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.95%  '
$ws.Range("E30").Value = '  -2.09%  '
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.33'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.93'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.84'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '165.97'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.39%  '
$ws.Range("E36").Value = '  -2.00%  '
$ws.Range("D37").Value = '3.402.76'
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("E38").Value = '  -3.52%  '
$ws.Range("E39").Value = '  -2.10%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.772'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.22%  '
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '25.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -9.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.32'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.47%  '
$ws.Range("E43").Value = '  -3.59%  '
$ws.Range("E44").Value = '  -2.48%  '
$ws.Range("D45").Value = '2.450.25'
$ws.Range("E45").Value = '  -3.73%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.60'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0258'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.81%  '
$ws.Range("E50").Value = '  -2.18%  '
$ws.Range("E51").Value = '  -3.45%  '
